# Update the "想去人数" (want-to-go count) figures in column F that were
# refreshed by the scraper run (output generated at 456a3b4).
#
# Sheet "展览": F2 1368->1374, F3 2902->2911, F4 11->13, F5 265->267
# Sheet "全部类型": F3 1368->1374, F4 2902->2911, F5 11->13, F7 265->267

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1374
$wsExhibit.Range("F3").Value = 2911
$wsExhibit.Range("F4").Value = 13
$wsExhibit.Range("F5").Value = 267

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1374
$wsAll.Range("F4").Value = 2911
$wsAll.Range("F5").Value = 13
$wsAll.Range("F7").Value = 267
